$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to Text format before assigning numeric-looking strings,
# so Excel does not auto-convert them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.514.26'
$ws.Range("E2").Value = '  +4.99%  '
$ws.Range("D3").Value = '3.491.84'
$ws.Range("E3").Value = '  +4.87%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '587.48'
$ws.Range("E5").Value = '  +6.13%  '
$ws.Range("D6").Value = '185.50'
$ws.Range("E6").Value = '  +6.98%  '
$ws.Range("E7").Value = '  +2.26%  '
$ws.Range("D8").Value = '3.483.18'
$ws.Range("E8").Value = '  +4.87%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  +2.33%  '
$ws.Range("D11").Value = '0.652'
$ws.Range("E11").Value = '  +3.26%  '
$ws.Range("D12").Value = '56.28'
$ws.Range("E12").Value = '  +5.44%  '
$ws.Range("D13").Value = '0.0000281'
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("D14").Value = '9.47'
$ws.Range("E14").Value = '  +4.05%  '
$ws.Range("D15").Value = '4.044.24'
$ws.Range("E15").Value = '  +4.88%  '
$ws.Range("D16").Value = '18.88'
$ws.Range("E16").Value = '  +4.23%  '
$ws.Range("D17").Value = '3.496.50'
$ws.Range("E17").Value = '  +4.40%  '
$ws.Range("D18").Value = '67.416.45'
$ws.Range("E18").Value = '  +5.08%  '
$ws.Range("D19").Value = '12.22'
$ws.Range("E19").Value = '  +4.32%  '
$ws.Range("E20").Value = '  -1.13%  '
$ws.Range("D22").Value = '491.09'
$ws.Range("E22").Value = '  +8.21%  '
$ws.Range("D23").Value = '5.44'
$ws.Range("E23").Value = '  +6.59%  '
$ws.Range("D24").Value = '16.59'
$ws.Range("E24").Value = '  +17.94%  '
$ws.Range("D25").Value = '4.47'
$ws.Range("E25").Value = '  +10.00%  '
$ws.Range("D26").Value = '90.13'
$ws.Range("E26").Value = '  +3.24%  '
$ws.Range("E27").Value = '  +3.13%  '
$ws.Range("D28").Value = '11.00'
$ws.Range("E28").Value = '  +3.68%  '
$ws.Range("D29").Value = '9.16'
$ws.Range("E29").Value = '  +6.72%  '
$ws.Range("D30").Value = '31.59'
$ws.Range("E30").Value = '  +2.08%  '
$ws.Range("D31").Value = '7.19'
$ws.Range("E31").Value = '  +10.48%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").Value = '11.75'
$ws.Range("E32").Value = '  +3.00%  '
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").Value = '64.46'
$ws.Range("E33").Value = '  +4.58%  '
$ws.Range("D34").Value = '595.88'
$ws.Range("E34").Value = '  +5.31%  '
$ws.Range("D35").Value = '0.112'
$ws.Range("E35").Value = '  +4.71%  '
$ws.Range("E36").Value = '  +6.44%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").Value = '36.67'
$ws.Range("E38").Value = '  +3.87%  '
$ws.Range("D39").Value = '3.57'
$ws.Range("E39").Value = '  +1.87%  '
$ws.Range("D40").Value = '0.387'
$ws.Range("E40").Value = '  +5.98%  '
$ws.Range("E41").Value = '  +5.62%  '
$ws.Range("D42").Value = '3.246.71'
$ws.Range("E42").Value = '  +6.13%  '
$ws.Range("D43").Value = '2.92'
$ws.Range("E43").Value = '  +6.30%  '
$ws.Range("E44").Value = '  +3.43%  '
$ws.Range("D45").Value = '2.54'
$ws.Range("E45").Value = '  +3.50%  '
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("D47").Value = '2.76'
$ws.Range("E47").Value = '  +22.32%  '
$ws.Range("E48").Value = '  +1.78%  '
$ws.Range("E49").Value = '  +12.35%  '
$ws.Range("D50").Value = '8.78'
$ws.Range("E50").Value = '  +7.57%  '
$ws.Range("E51").Value = '  -0.07%  '

# Restore original (no explicit style) formatting on column D.
$ws.Range("D2:D51").ClearFormats()
